$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "along" values live in column A; find the last used row, then scan from
# bottom to top (so deleting a row doesn't disturb the indices of rows we
# still need to examine) and remove every row whose "along" value is
# "Saudi Arabia". This drops Saudi Arabia from each of the three country
# groups and re-runs/renumbers the remaining rows, matching the updated
# (cropped) dataset.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = $lastRow; $r -ge 2; $r--) {
    $val = $ws.Cells.Item($r, 1).Value()
    if ($val -eq "Saudi Arabia") {
        $ws.Rows.Item($r).Delete()
    }
}
